# ---------------------------------------------------------------------------
# Add a new "2022-Q1" sheet (before the "总计" summary sheet) and add a
# corresponding summary row to "总计".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" worksheet, positioned right before "总计".
#    We copy the layout (and therefore formatting/styles) of an existing
#    quarterly sheet ("2021-Q4") and then overwrite the values.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$beforeSheet = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q1"

$template.Range("A1:H9").Copy($newSheet.Range("A1"))

# The header row has no value in column A (only the data rows do) - drop
# the stray cell that Copy left behind.
$newSheet.Range("A1").ClearContents()

# Rows 10 and 11 don't exist on the 9-row template - give them the same
# formatting as the other data rows before filling in their values.
$template.Range("A9:H9").Copy()
$newSheet.Range("A10:H10").PasteSpecial(-4122)
$newSheet.Range("A11:H11").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows
$rows = @(
    @(0, "512980", "广发中证传媒ETF",               "44.11", "99.38", "3.05", "1.3454", 10),
    @(1, "159869", "华夏中证动漫游戏ETF",            "6.20",  "98.75", "6.27", "0.3887", 6),
    @(2, "516010", "国泰中证动漫游戏ETF",            "4.95",  "98.91", "6.15", "0.3044", 6),
    @(3, "001628", "招商体育文化休闲股票",            "2.95",  "83.21", "4.11", "0.1212", 7),
    @(4, "161030", "富国中证体育产业指数",            "2.32",  "93.75", "5.09", "0.1181", 2),
    @(5, "516770", "华泰柏瑞中证动漫游戏ETF",         "1.11",  "96.56", "6.12", "0.0679", 6),
    @(6, "164818", "工银瑞信中证传媒指数（LOF）A",     "1.99",  "92.70", "2.83", "0.0563", 10),
    @(7, "159805", "鹏华中证传媒ETF",                "1.73",  "96.29", "2.97", "0.0514", 10),
    @(8, "010677", "工银瑞信中证传媒指数（LOF）C",     "0.25",  "92.70", "2.83", "0.0071", 10),
    @(9, "165522", "信诚中证TMT产业主题指数（LOF）",   "0.58",  "93.74", "1.20", "0.0070", 4)
)

foreach ($row in $rows) {
    $r = 2 + $row[0]
    # Columns B-G hold text values in the source data (fund code, name,
    # scale, position, etc. are all stored as strings, not numbers). A
    # leading apostrophe forces Excel to keep them as text instead of
    # auto-converting to numbers (and losing leading zeros, for example).
    $newSheet.Range("A$r").Value = $row[0]
    $newSheet.Range("B$r").Value = "'" + $row[1]
    $newSheet.Range("C$r").Value = "'" + $row[2]
    $newSheet.Range("D$r").Value = "'" + $row[3]
    $newSheet.Range("E$r").Value = "'" + $row[4]
    $newSheet.Range("F$r").Value = "'" + $row[5]
    $newSheet.Range("G$r").Value = "'" + $row[6]
    $newSheet.Range("H$r").Value = $row[7]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new summary row for 2022-Q1 above the
#    existing rows (the newest quarter is always listed first).
#    NOTE: re-fetch the sheet by name (rather than reusing $beforeSheet)
#    since the previous reference now tracks the newly inserted/renamed
#    sheet.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Shift the existing data rows (2-6) down by one (working bottom-up so we
# never clobber a row before reading it).
for ($r = 6; $r -ge 2; $r--) {
    $src = $r
    $dst = $r + 1
    $av = $total.Range("A$src").Value2
    $bv = $total.Range("B$src").Value2
    $cv = $total.Range("C$src").Value2
    $dv = $total.Range("D$src").Value2
    $total.Range("A$dst").Value = $av
    $total.Range("B$dst").Value = $bv
    $total.Range("C$dst").Value = $cv
    $total.Range("D$dst").Value = $dv
}

# Row 7 is brand new - give it the same formatting as row 6 before writing.
$total.Range("A6:D6").Copy()
$total.Range("A7:D7").PasteSpecial(-4122)

# Fill in the new 2022-Q1 summary row.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 10
$total.Range("D2").Value = 2.47

# Renumber the index column (A) for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
